$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '67.560.94'
$ws.Range("E2").Value2 = '  -1.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '3.513.11'
$ws.Range("E3").Value2 = '  -4.10%  '
$ws.Range("E4").Value2 = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '198.04'
$ws.Range("E5").Value2 = '  -2.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '552.69'
$ws.Range("E6").Value2 = '  -4.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.631'
$ws.Range("E7").Value2 = '  +1.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '3.503.08'
$ws.Range("E8").Value2 = '  -4.25%  '
$ws.Range("E9").Value2 = '  -0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '0.654'
$ws.Range("E10").Value2 = '  -4.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '61.11'
$ws.Range("E11").Value2 = '  +6.55%  '
$ws.Range("E12").Value2 = '  -8.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '0.0000268'
$ws.Range("E13").Value2 = '  -9.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '9.79'
$ws.Range("E14").Value2 = '  -3.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '4.071.64'
$ws.Range("E15").Value2 = '  -4.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '3.514.88'
$ws.Range("E16").Value2 = '  -4.27%  '
$ws.Range("E17").Value2 = '  -1.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '67.266.61'
$ws.Range("E18").Value2 = '  -2.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '18.33'
$ws.Range("E19").Value2 = '  -2.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '11.82'
$ws.Range("E20").Value2 = '  -6.23%  '
$ws.Range("E21").Value2 = '  -6.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '393.70'
$ws.Range("E22").Value2 = '  -2.88%  '
$ws.Range("B23").Value2 = 'RenderToken'
$ws.Range("C23").Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '11.93'
$ws.Range("E23").Value2 = '  -8.26%  '
$ws.Range("B24").Value2 = 'PancakeSwap'
$ws.Range("C24").Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '3.96'
$ws.Range("E24").Value2 = '  -7.32%  '
$ws.Range("B25").Value2 = 'Litecoin'
$ws.Range("C25").Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '85.19'
$ws.Range("E25").Value2 = '  -1.39%  '
$ws.Range("E26").Value2 = '  -0.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '12.29'
$ws.Range("E27").Value2 = '  -3.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '2.81'
$ws.Range("E28").Value2 = '  -5.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '8.86'
$ws.Range("E29").Value2 = '  -4.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '717.47'
$ws.Range("E30").Value2 = '  +2.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '31.27'
$ws.Range("E31").Value2 = '  -2.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '7.00'
$ws.Range("E32").Value2 = '  -15.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '11.69'
$ws.Range("E33").Value2 = '  -5.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '63.98'
$ws.Range("E34").Value2 = '  -1.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '0.111'
$ws.Range("E35").Value2 = '  -5.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '38.36'
$ws.Range("E36").Value2 = '  -10.76%  '
$ws.Range("E37").Value2 = '  -0.03%  '
$ws.Range("E38").Value2 = '  -9.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '3.01'
$ws.Range("E39").Value2 = '  -5.02%  '
$ws.Range("B40").Value2 = 'FirstDigitalUSD'
$ws.Range("C40").Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '0.999'
$ws.Range("E40").Value2 = '  -0.17%  '
$ws.Range("B41").Value2 = 'Kaspa'
$ws.Range("C41").Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '0.131'
$ws.Range("E41").Value2 = '  -8.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '3.059.26'
$ws.Range("E42").Value2 = '  -5.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '0.0₃0679'
$ws.Range("E43").Value2 = '  -15.78%  '
$ws.Range("E44").Value2 = '  +4.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '2.51'
$ws.Range("E45").Value2 = '  -12.15%  '
$ws.Range("E46").Value2 = '  -4.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '0.130'
$ws.Range("E47").Value2 = '  -1.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '2.55'
$ws.Range("E48").Value2 = '  -15.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '138.67'
$ws.Range("E49").Value2 = '  -2.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '2.93'
$ws.Range("E50").Value2 = '  -5.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '8.20'
$ws.Range("E51").Value2 = '  -8.85%  '
